# Optuna Attempt (go back with original)
# Update forecast figures on the "Forecast Comparison" sheet and the
# derived totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- "Forecast Comparison" sheet -------------------------------------

# MyForecast (column D), rows 3-14: 7 -> 6
for ($r = 3; $r -le 14; $r++) {
    $wsForecast.Range("D$r").Value = 6
}

# Inventory Coverage (column H), rows 2-17
$hValues = @{
    2  = 18.46
    3  = 17.46
    4  = 16.46
    5  = 15.46
    6  = 14.92
    7  = 13.49
    8  = 12.49
    9  = 11.49
    10 = 10.49
    11 = 9.49
    12 = 8.49
    13 = 7.49
    14 = 6.49
    15 = 5.36
    16 = 4.47
    17 = 3.47
}
foreach ($r in $hValues.Keys) {
    $wsForecast.Range("H$r").Value = $hValues[$r]
}

# Seasonality Index (column L), rows 2-17
$lValues = @{
    2  = 1.08
    3  = 0.98
    4  = 1.05
    5  = 1.17
    6  = 0.8100000000000001
    7  = 1.1
    8  = 1.11
    9  = 1.04
    10 = 1.14
    11 = 1.04
    12 = 0.89
    13 = 1.05
    14 = 1
    15 = 1.18
    16 = 0.92
    17 = 0.87
}
foreach ($r in $lValues.Keys) {
    $wsForecast.Range("L$r").Value = $lValues[$r]
}

# --- "Summary" sheet ---------------------------------------------------
# These values are stored as text (not numbers) in the workbook, so a
# leading apostrophe forces Excel to keep them as plain text instead of
# auto-converting the numeric-looking string into a number.

$wsSummary.Range("B9").Value  = "'104"
$wsSummary.Range("B10").Value = "'52"
$wsSummary.Range("B11").Value = "'26"
$wsSummary.Range("B12").Value = "'7"
$wsSummary.Range("B14").Value = "'6"
